# Fill English translations (column C) for the "14B" sheet (Wesnoth Marbuss Escape dialogues)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("14B")
$ws.Activate()

$ws.Cells.Item(70, 3).Value = 'We meet again, Marbus. This time on opposite sides.'
$ws.Cells.Item(71, 3).Value = 'Jarl Oferiu! I didn''t think we''ll meet again! You did a great job helping us defead the Wild Gon. I could say, we wouldn''t be here without you! You must be proud, aren''t you?'
$ws.Cells.Item(72, 3).Value = 'Silence, Gewold. I''m glad to see you in good health, Jarl. Your help was really much for us. Is there anything I can do for you?'
$ws.Cells.Item(73, 3).Value = 'You can gather the Great Horde and leave the North, but I think it''s too much of a demand? Well… Then I don''t want anything from you. I''ll stand my ground and protect lady Beatrice till I die.'
$ws.Cells.Item(74, 3).Value = 'I respect that, Jarl. I think I could leave you alive and make a supervisor of slaves… An administrator of your own race… But I''m afraid you''re too proud for that. You wouldn''t take that job, or rather, commit suicide instead. Anyway, we''re thankful for your assistance against Wild Gon. We promise to kill you quickly. Painlessly.'
$ws.Cells.Item(75, 3).Value = 'Stop here, damned orc! I swear on my sword, you won''t pass any further!'
$ws.Cells.Item(76, 3).Value = 'Stop. Halt. Go away. Those words work on me just like "cursed" and "filthy". It''s rather boring… You are boring, Jarl. Even standing before the Great Horde, you cannot say any better words. But I see, you aren''t an elder… Young, wet behind the ears… I''ll give you a second chance. Go, say something elevated!'
$ws.Cells.Item(77, 3).Value = 'No? Well, it''s your choice. Gewold, gut him before his warrior''s eyes.'
$ws.Cells.Item(78, 3).Value = 'You don''t have to say it twice, Great Sovereign!'
$ws.Cells.Item(79, 3).Value = 'Ugh… Dammit… But at least I die with honour… I did everything I could…'
$ws.Cells.Item(80, 3).Value = 'You weren''t really able to do much, stupid girl… But when it comes to me, I can do a lot. And I always keep my word. I''ll remove your name from history.'
$ws.Cells.Item(81, 3).Value = 'Dammit! Betray my own race for future glory… And die on the same day… Marbus… You didn''t keep your promise…'
$ws.Cells.Item(82, 3).Value = 'Gewold! You were supposed to keep an eye on her! I swear, when this battle will be over, you''ll pay for it!'
$ws.Cells.Item(83, 3).Value = 'She way an adult… She was a warrior. I came here to kill, not to take care of your girls!'
$ws.Cells.Item(84, 3).Value = 'This Person thinks, that you should watch your words… You won''t go without a punishment, but we''ll talk about it after the battle.'
$ws.Cells.Item(85, 3).Value = 'Ugh… Dammit… I can''t fight anymore…'
$ws.Cells.Item(86, 3).Value = 'Should we gut him?'
$ws.Cells.Item(87, 3).Value = 'No. As I said, we owe you a lot, Jarl. So you''ll die with honour. Lightest death in the entire battle. Give him a dagger.'
$ws.Cells.Item(88, 3).Value = 'Be damned, Marbus… You''ll fall one day. One day, the Great Horde will fall, same as when the Northern Alliance was formed… History will repeat itself… Orcs cannot rule over us…'
$ws.Cells.Item(89, 3).Value = 'I''m waiting, Jarl.'
$ws.Cells.Item(90, 3).Value = '*Slits his wrists*'
$ws.Cells.Item(91, 3).Value = 'You''ll… Fall….'
$ws.Cells.Item(92, 3).Value = 'Ugh… What… What are you doing...?'
$ws.Cells.Item(93, 3).Value = 'You didn''t hear the Great Sovereign? We''re gutting you alive!'
$ws.Cells.Item(94, 3).Value = 'Aaaa! No, please… Ugh…'
$ws.Cells.Item(95, 3).Value = 'This is how fools die…'
$ws.Cells.Item(96, 3).Value = 'It took you a while to reach my defense line, orcs. I couldnt wait to see you. Same, as my friends, who you know well.'
$ws.Cells.Item(97, 3).Value = 'Gilfit the Snowfeet'
$ws.Cells.Item(101, 3).Value = 'The Snowfeet Clan will fight you to the last soldier! I''m Gilfit, last of my kin. You rended through our caves to flee from undeads. Even now, when you are feared in the whole North, we remember you fleeing for your life.'
$ws.Cells.Item(102, 3).Value = 'Your clan… This Person hates it… I was hiding in shadows of your waters, until I met the Great Sovereign… Today it''s time for revenge. Armors and bones of dwarves are hard… It''s good then, that swords of This Person are also made from very good steel.'
$ws.Cells.Item(103, 3).Value = 'Dwarves? And here I thought that we exterminated all pests underground. It''s good they showed themselves on the surface. We''ll feast our eyes on those dwarfs until we finally finish off this race and make them our slaves. And you, Jarl… You''re old, embittered, and full of hatred. It''s time to end your life.'
$ws.Cells.Item(104, 3).Value = 'I remember more battles than years you lived. I surpass you in tactical genius, experience, knowledge… You won''t be able to overcome my soldiers.'
$ws.Cells.Item(105, 3).Value = 'We will, because we surpass you in numbers, bloodlust, and equipment. You don''t even know how many weapons we''ve got… We keep most of it in storages, because even the Great Horde can''t use all of it. And you? It''s probably rather hard to get any metal in this forest. How do you fight? You share one sword for three wariors?'
$ws.Cells.Item(108, 3).Value = 'No… Without me… The Snowfeet Clan…'
$ws.Cells.Item(106, 3).Value = 'Even if we lack equipment, we still have the spirit! Elves, dwarves! Let''s pay them back for everything they did to us!'
$ws.Cells.Item(109, 3).Value = 'Your clan is the past now. Same goes for all the Knalga. We''ll stumble on you, enslave and make sure, that all your pride and hubris will be gone, replaced by fear for the Great Horde. None of underground''s sons will dare to raise an axe against his orcish masters!'
$ws.Cells.Item(110, 3).Value = 'Curse you… To hell… Dwarves won''t… Ever… Bow to you…'
$ws.Cells.Item(111, 3).Value = 'Your pieces will bow to me! Tear him into shreds and throw the remains before my feet!'
$ws.Cells.Item(112, 3).Value = 'So many years of experience… Fights… Noble battles and duels… To be attacked by… Filthy orcs! You should all die when the Wild Gon attacked. You didn''t do, though… Instead, you spread, like… like vermin! Like cockroaches or ants! You spread your monstrosity... You destroy everything that''s beautiful, worth fighting for...'
$ws.Cells.Item(113, 3).Value = 'Well, many people shares similar opinion about us. The problem is, we don''t really care about it. Same goes for talking with you… Farewell, Jarl… I wish you won''t be so bitter old greybeard in the afterlife.'
$ws.Cells.Item(124, 3).Value = 'So that''s how death looks like… It''s so cold… And lonely…'
$ws.Cells.Item(125, 3).Value = 'Any death is better than what we''re preparing for those who survive. Living under the Great Horde''s rule will be a greater suffering for them than any physical pain. And each, and every suffering one will be able to blame it on you. It''s you, who lost. You let them down... Die with this thought, Beatrice.'
$ws.Cells.Item(126, 3).Value = 'Finish it…'
$ws.Cells.Item(127, 3).Value = 'Farewell, elvish witch. This is how your story ends and a brand new story starts. A story, about invincible Great Horde, ruling the whole North!'
$ws.Cells.Item(114, 3).Value = 'And here we have an elvish witch, a descendant of Kalenz, Beatrice. I was dreaming about that moment, when I can oppose you for years. For all those moments when you openly despised me, treated like I was worse than you… You don''t know how much I hate you...'
$ws.Cells.Item(115, 3).Value = 'Your betrayal is more painful for me than deaths of all those who protected this place. They died in glory, fighting for what they believed. And you? You have choosen to live with the Great Horde, with those, who murdered your kin. You latched onto that grace and betrayed us.'
$ws.Cells.Item(116, 3).Value = 'I did what''s the best for me and the Blood Elves. We argued for years, it''s time to settle this conflict. Rise above, have you under my feet and look on you, living as slaves - that''s what I want!'
$ws.Cells.Item(117, 3).Value = 'You become as despicable as orcs, girl. But it was only predictable… You were always full of hatred, unhealthy fascination about killing and lust for power. How could you think you could ever be equal to us? I tell you, even if my people will be forced to clean your shoes, everyone will still despise you inside.'
$ws.Cells.Item(118, 3).Value = 'Marbus… I regret meeting you again in such circumstances. I don''t know what has come into you. Why you let your race become posessed by bloodlust again? Why you let this situation happen? Don''t you see how everyone suffer because of your actions?'
$ws.Cells.Item(119, 3).Value = 'My warriors don''t suffer. I lead the Great Horde through a bright path, that it wants to be lead. And when it comes to reasons… Well, everything could be different. But, provided my cousin lived. If you really tried to keep your alliance with us.'
$ws.Cells.Item(120, 3).Value = 'About your cousin I only known after it was done… I didn''t want his death. I believed entirely that you can change your destiny. But it didn''t happen. You still do everything to be hated. After a thousand years the King of Gon will return. Before I thought, that he''s our curse, but now I understand - he was our chance. A chance to get rid of you. With a great price, but still a chance...'
$ws.Cells.Item(121, 3).Value = 'It doesn''t matter anymore. The King of Gon is dead, and even if he returns, he''ll find us even more powerfull than in this millenium. The Gjallarhorn will be passed to every new Great Sovereign, on every generation… We''ll be ready for this curse to come.'
$ws.Cells.Item(122, 3).Value = 'Talking won''t bring us anything more. Attack… I''m ready for the strike…'
$ws.Cells.Item(123, 3).Value = 'And I''ll be glad to deal it. Forward, my warriors! Surround her fortress, don''t let any tree-hugger escape! It''s time for our final triumph! Tash''kfug adin!'
$ws.Cells.Item(63, 3).Value = '*Gets down on one knee and bows her head. Her fighters lower their bows and do likewise.*'

# Row heights grew for two rows whose new English text wraps onto more lines
# than the previous (empty / shorter) content did.
$ws.Rows.Item(84).RowHeight = 38.25
$ws.Rows.Item(124).RowHeight = 25.5

# Restore the view roughly where the translator left off (scrolled down, cell L61 selected)
$excel.ActiveWindow.ScrollRow = 104
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L61").Select()
